$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.214.34'
$ws.Range("E2").Value = '  +6.82%  '
$ws.Range("D3").Value = '3.676.67'
$ws.Range("E3").Value = '  +19.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.53'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.18%  '
$ws.Range("D7").Value = '3.673.78'
$ws.Range("E7").Value = '  +19.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.63%  '
$ws.Range("E10").Value = '  +7.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.57'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.44%  '
$ws.Range("E12").Value = '  +6.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.94'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +11.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000254'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.47%  '
$ws.Range("D15").Value = '4.278.05'
$ws.Range("E15").Value = '  +18.80%  '
$ws.Range("D16").Value = '71.145.19'
$ws.Range("E16").Value = '  +6.78%  '
$ws.Range("D17").Value = '3.665.04'
$ws.Range("E17").Value = '  +18.68%  '
$ws.Range("E18").Value = '  +2.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '519.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +17.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.745'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.89%  '
$ws.Range("E26").Value = '  +8.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.38%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  +12.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +14.08%  '
$ws.Range("E32").Value = '  +6.92%  '
$ws.Range("E33").Value = '  +17.63%  '
$ws.Range("E34").Value = '  +4.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.20'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.01'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.343'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.29%  '
$ws.Range("E39").Value = '  +8.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '46.32'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.64%  '
$ws.Range("E42").Value = '  +4.07%  '
$ws.Range("D43").Value = '3.188.18'
$ws.Range("E43").Value = '  +14.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '400.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.20%  '
$ws.Range("E47").Value = '  +6.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '28.23'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +15.57%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '135.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.44'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +12.32%  '
